$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 71, pushing existing rows 71-145 down to 72-146
$ws.Rows.Item(71).Insert()

# Populate the new row 71 with the latest weekly record
$ws.Range("A71").Value = 10
$ws.Range("B71").Value = "Vega Modelo de Temuco"
$ws.Range("C71").Value = "La Araucanía"
$ws.Range("D71").Value = 44483
$ws.Range("E71").Value = 9
$ws.Range("F71").Value = 100112052
$ws.Range("G71").Value = "Albahaca"
$ws.Range("H71").Value = "Sin especificar"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 30
$ws.Range("K71").Value = 7000
$ws.Range("L71").Value = 7000
$ws.Range("M71").Value = 7000
$ws.Range("N71").Value = "$/paquete"
$ws.Range("O71").Value = "Región de Arica y Parinacota"
$ws.Range("P71").Value = 7000
$ws.Range("Q71").Value = 1
$ws.Range("R71").Value = "Hortaliza"
